# Atualizado por script em 01-12-2023 20:45
# Adds two new match rows (144 and 145) to the Ekstraklasa 2023-2024 sheet,
# mirroring the data/formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-MatchRow($Row, $Indice, $Pais, $Torneio, $Temporada, $DataPartida, $Home, $HomeGols, $Away, $AwayGols, $HomeOpenOdds, $HomeOpenData, $HomeCloseOdds, $HomeCloseData, $DrawOpenOdds, $DrawOpenData, $DrawCloseOdds, $DrawCloseData, $AwayOpenOdds, $AwayOpenData, $AwayCloseOdds, $AwayCloseData, $Url) {
    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = $Pais
    $ws.Cells.Item($Row, 3).Value = $Torneio
    $ws.Cells.Item($Row, 4).Value = $Temporada
    $ws.Cells.Item($Row, 5).Value = $DataPartida
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGols
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGols
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenData
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseData
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenData
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseData
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenData
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseData
    $ws.Cells.Item($Row, 22).Value = $Url

    # Match the existing styling: column A (bold/boxed/centered) and
    # column E (datetime number format) are copied from the previous
    # data row so they reuse the same cell-style indexes as the rest
    # of the sheet instead of generating brand-new styles.
    $prevRow = $Row - 1
    $ws.Range("A" + $prevRow).Copy() | Out-Null
    $ws.Range("A" + $Row).PasteSpecial(-4122) | Out-Null
    $ws.Range("E" + $prevRow).Copy() | Out-Null
    $ws.Range("E" + $Row).PasteSpecial(-4122) | Out-Null
}

Add-MatchRow 144 143 "poland" "ekstraklasa" "2023-2024" 45261.75 "Warta Poznan" 1 "Jagiellonia" 2 2.68 "25/11/2023 20:13" 3.52 "01/12/2023 17:58" 3.11 "25/11/2023 20:13" 3.17 "01/12/2023 17:59" 2.88 "25/11/2023 20:13" 2.3 "01/12/2023 17:59" "https://www.betexplorer.com/football/poland/ekstraklasa/warta-poznan-jagiellonia/M3pirQS0/"

Add-MatchRow 145 144 "poland" "ekstraklasa" "2023-2024" 45261.85416666666 "Gornik Zabrze" 1 "Pogon Szczecin" 0 3.59 "26/11/2023 11:42" 3.15 "01/12/2023 20:24" 3.44 "26/11/2023 11:42" 3.82 "01/12/2023 20:28" 2.03 "26/11/2023 11:42" 2.2 "01/12/2023 20:28" "https://www.betexplorer.com/football/poland/ekstraklasa/gornik-zabrze-pogon-szczecin/IRV63za9/"
